{"js": "// 1. Remove the _GoBack bookmark from the title paragraph (it will be\n//    re-added later, anchored to the new last-but-one paragraph).\nconst existingGoBack = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nexistingGoBack.load(\"isNullObject\");\nawait context.sync();\nif (!existingGoBack.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 2. Merge the split \">>\", \">  your\", \" stuff after this line >>>\" runs\n//    (with proofErr markers in between) into a single run of text\n//    \">>>  your stuff after this line >>>\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nconst mergedText = \">>>  your stuff after this line >>>\";\nconst newParaText =\n  \"This is Samandeep, trying to play around with GitHub and see what all commands do we have.\";\n\nlet mergeTarget = null;\nlet replaceTarget = null;\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === mergedText) {\n    mergeTarget = items[i];\n  }\n}\n// 3. The paragraph to receive the new sentence is the second-to-last\n//    paragraph (an empty <w:p/> right before the final empty <w:p/>).\nreplaceTarget = items[items.length - 2];\n\nif (mergeTarget) {\n  mergeTarget.clear();\n  await context.sync();\n  mergeTarget.insertText(mergedText, Word.InsertLocation.start);\n  await context.sync();\n}\n\nreplaceTarget.insertText(newParaText, Word.InsertLocation.start);\nawait context.sync();\n\nconst endRange = replaceTarget.getRange(Word.RangeLocation.end);\nendRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Remove the _GoBack bookmark from the title paragraph (it will be\n#    re-added later in the new last-but-one paragraph).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2. Merge the split \">>\", \">  your\", \" stuff after this line >>>\" runs\n#    (with proofErr markers in between) into a single run of text\n#    \">>>  your stuff after this line >>>\".\n$target = \">>>  your stuff after this line >>>\"\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    $rng = $para.Range\n    [void]$rng.MoveEnd(1, -1)\n    if ($rng.Text -eq $target) {\n        $rng.Delete()\n        $rng.InsertAfter($target)\n        break\n    }\n}\n\n# 3. Replace the empty paragraph right before the final empty paragraph\n#    with the new sentence, followed by a fresh (zero-length) _GoBack\n#    bookmark positioned right after the inserted text. We add the\n#    bookmark first (while the paragraph is still empty) and then type\n#    the sentence in *before* the bookmark, so the bookmark ends up\n#    sitting right after the text instead of wrapping it.\n$count = $d.Paragraphs.Count\n$newPara = $d.Paragraphs.Item($count - 1)\n$newRng = $newPara.Range\n[void]$newRng.MoveEnd(1, -1)\n\n$d.Bookmarks.Add(\"_GoBack\", $newRng)\n$bmRng = $d.Bookmarks.Item(\"_GoBack\").Range\n$bmRng.InsertBefore(\"This is Samandeep, trying to play around with GitHub and see what all commands do we have.\")\n"}
